# Remove the <w:contextualSpacing w:val="0"/> element from every
# paragraph's paragraph properties (w:pPr). This mirrors the commit's
# XML diff, which deletes that single child element from <w:pPr> in
# every paragraph that had it, leaving everything else untouched.

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pXml = $p.XML()

    if ($pXml -like "*<w:contextualSpacing*") {
        $newXml = $pXml -replace '<w:contextualSpacing[^>]*/>', ''
        [void]$p.Range.InsertXML($newXml)
    }
}
